$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trend-check columns (I/J) computed from the most recent 5 years of
# quarterly wetland extent data (rows 6-13), comparing each 2010 quarter
# to the corresponding 2005 quarter.
$ws.Range("I13").Formula = "=1-(C10/C6)"
$ws.Range("J13").Formula = "=-I13/5"

$ws.Range("I14").Formula = "=1-(C11/C7)"
$ws.Range("J14").Formula = "=-I14/5"

$ws.Range("I15").Formula = "=1-(C12/C8)"
$ws.Range("J15").Formula = "=-I15/5"

$ws.Range("I16").Formula = "=1-(C13/C9)"
$ws.Range("J16").Formula = "=-I16/5"

# Update the selected range to reflect the newly added trend cells.
$ws.Range("J13:J16").Select()
